# Atualização da base de dados
# Applies updated enrollment figures to the "Resumo Inscrições" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each entry: Row, Column -> New Value
$updates = @(
    @{ Row = 4;  Col = "E"; Value = 5 },

    @{ Row = 15; Col = "E"; Value = 86 },
    @{ Row = 15; Col = "F"; Value = 37 },
    @{ Row = 15; Col = "H"; Value = 37 },

    @{ Row = 17; Col = "E"; Value = 51 },

    @{ Row = 18; Col = "E"; Value = 46 },

    @{ Row = 19; Col = "E"; Value = 24 },

    @{ Row = 33; Col = "E"; Value = 15 },

    @{ Row = 41; Col = "E"; Value = 15 },
    @{ Row = 41; Col = "F"; Value = 6 },
    @{ Row = 41; Col = "H"; Value = 6 },

    @{ Row = 46; Col = "E"; Value = 13 },
    @{ Row = 46; Col = "F"; Value = 1 },
    @{ Row = 46; Col = "H"; Value = 1 },

    @{ Row = 48; Col = "E"; Value = 9 },
    @{ Row = 48; Col = "F"; Value = 5 },
    @{ Row = 48; Col = "H"; Value = 5 },

    @{ Row = 49; Col = "E"; Value = 25 },
    @{ Row = 49; Col = "F"; Value = 12 },
    @{ Row = 49; Col = "H"; Value = 12 },

    @{ Row = 64; Col = "E"; Value = 20 },
    @{ Row = 64; Col = "F"; Value = 10 },
    @{ Row = 64; Col = "H"; Value = 10 },

    @{ Row = 65; Col = "E"; Value = 14 },

    @{ Row = 66; Col = "E"; Value = 18 },

    @{ Row = 71; Col = "E"; Value = 9 },

    @{ Row = 79; Col = "E"; Value = 10 },

    @{ Row = 80; Col = "E"; Value = 15 }
)

foreach ($u in $updates) {
    $addr = "$($u.Col)$($u.Row)"
    $ws.Range($addr).Value = $u.Value
}
